$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Cells whose new text would otherwise be auto-parsed as a plain number
# (e.g. "216.73") are first forced to Text format so they stay strings,
# matching the workbook's original inline-string cell content.
$ws.Range("D2").Value = "27.146.62"
$ws.Range("E2").Value = "  +1.00%  "
$ws.Range("D3").Value = "1.639.64"
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.73"
$ws.Range("E5").Value = "  -0.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.522"
$ws.Range("E6").Value = "  +2.26%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.255"
$ws.Range("E8").Value = "  -0.16%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.0626"
$ws.Range("E9").Value = "  +0.07%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.97"
$ws.Range("E10").Value = "  +0.54%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0846"
$ws.Range("E11").Value = "  -0.19%  "
$ws.Range("D12").Value = "1.867.37"
$ws.Range("E12").Value = "  -0.06%  "
$ws.Range("D13").Value = "1.641.70"
$ws.Range("E13").Value = "  +0.05%  "
$ws.Range("E14").Value = "  +0.29%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.542"
$ws.Range("E15").Value = "  +2.22%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "66.84"
$ws.Range("E16").Value = "  -0.56%  "
$ws.Range("D17").Value = "27.136.82"
$ws.Range("D18").Value = "0.0₃0740"
$ws.Range("E18").Value = "  +1.22%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "217.13"
$ws.Range("E19").Value = "  -1.11%  "
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.96"
$ws.Range("E21").Value = "  +1.59%  "
$ws.Range("E22").Value = "  +3.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "9.13"
$ws.Range("E24").Value = "  -0.32%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "147.13"
$ws.Range("E25").Value = "  +0.11%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.43"
$ws.Range("E27").Value = "  +1.05%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.68"
$ws.Range("E29").Value = "  -0.68%  "
$ws.Range("E30").Value = "  +0.79%  "
$ws.Range("E31").Value = "  +0.03%  "
$ws.Range("E32").Value = "  +1.52%  "
$ws.Range("E33").Value = "  +0.66%  "
$ws.Range("D34").Value = "1.307.15"
$ws.Range("E34").Value = "  +2.96%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.57"
$ws.Range("E35").Value = "  -0.09%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.47"
$ws.Range("E36").Value = "  +1.59%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0176"
$ws.Range("E37").Value = "  -1.35%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.859"
$ws.Range("E38").Value = "  +2.92%  "
$ws.Range("E39").Value = "  +1.86%  "
$ws.Range("E40").Value = "  +0.07%  "
$ws.Range("E41").Value = "  +0.01%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.22"
$ws.Range("E42").Value = "  +5.47%  "
$ws.Range("E43").Value = "  -1.92%  "
$ws.Range("D44").Value = "1.778.68"
$ws.Range("E44").Value = "  -0.04%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "61.74"
$ws.Range("E45").Value = "  -0.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "91.52"
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("E48").Value = "  +1.70%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0512"
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.67"
$ws.Range("E50").Value = "  +0.34%  "
$ws.Range("E51").Value = "  -0.21%  "
